# Auto-generated edit script applying the cryptos.xlsx price-table update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.739.00'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.850.57'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.00'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4327'
$ws.Range('E7').Value = '  +1.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3659'
$ws.Range('E8').Value = '  -0.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.05'
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07340'
$ws.Range('E10').Value = '  +1.24%  '
$ws.Range('E11').Value = '  -2.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.74'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.805.57'
$ws.Range('E13').Value = '  -3.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.347'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.533'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06935'
$ws.Range('E16').Value = '  +1.32%  '
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '80.36'
$ws.Range('E18').Value = '  +3.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000009046'
$ws.Range('E19').Value = '  +2.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.003'
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.41'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '27.867.37'
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.984'
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('E24').Value = '  -2.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.124.18'
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.988'
$ws.Range('E26').Value = '  -2.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '155.66'
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.62'
$ws.Range('E28').Value = '  +2.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '120.89'
$ws.Range('E29').Value = '  +8.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.271'
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.869'
$ws.Range('E31').Value = '  +1.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08917'
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7582'
$ws.Range('E33').Value = '  -1.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.544'
$ws.Range('E34').Value = '  -0.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.953'
$ws.Range('E35').Value = '  +1.24%  '
$ws.Range('E36').Value = '  +3.42%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.001'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.106'
$ws.Range('E38').Value = '  +0.97%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05420'
$ws.Range('E39').Value = '  +0.67%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01936'
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.833'
$ws.Range('E41').Value = '  -3.91%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5098'
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1661'
$ws.Range('E43').Value = '  +1.22%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.668'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.336'
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.45'
$ws.Range('E46').Value = '  +0.98%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.06539'
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4666'
$ws.Range('E48').Value = '  -1.24%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.45'
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.001'
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.624'
$ws.Range('E51').Value = '  -0.94%  '
